$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Perfil" backlog item (row 9 of the sheet / 7th data row of the Excel
# table "Tabela1") was removed. Deleting the whole sheet row shifts every
# row below it up by one, shrinks the table's range automatically, and
# drops the now-unused shared strings ("Perfil" and its description).
$ws.Rows(9).Delete()

# Match the post-edit selection (the user's cursor ended up on the row that
# shifted into the old "Perfil" row's place).
$ws.Range("C9:F9").Select()
